# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" quarterly sheet (holding the latest fund data)
# right after the "总计" summary sheet, pushing the existing "2022-Q2"
# and "2021-Q4" sheets one position to the right, and updates the "总计"
# summary sheet with the new quarter's totals.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Duplicate the existing "2022-Q2" sheet (so we inherit its layout/styles)
# and place the copy right after "总计"; this becomes the new "2022-Q4" sheet.
$q2Sheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Worksheets.Item($totalSheet.Index + 1)
$q4Sheet.Name = "2022-Q4"

# --- Populate the new "2022-Q4" sheet with the latest fund holdings ---
# Row 2: 010447 中邮未来成长混合A
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'010447"
$q4Sheet.Range("C2").Value = "中邮未来成长混合A"
$q4Sheet.Range("D2").Value = "'0.48"
$q4Sheet.Range("E2").Value = "'91.28"
$q4Sheet.Range("F2").Value = "'2.56"
$q4Sheet.Range("G2").Value = "'0.0123"
$q4Sheet.Range("H2").Value = 10

# Row 3: 010448 中邮未来成长混合C
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "'010448"
$q4Sheet.Range("C3").Value = "中邮未来成长混合C"
$q4Sheet.Range("D3").Value = "'0.07"
$q4Sheet.Range("E3").Value = "'91.28"
$q4Sheet.Range("F3").Value = "'2.56"
$q4Sheet.Range("G3").Value = "'0.0018"
$q4Sheet.Range("H3").Value = 10

# --- Update the "总计" summary sheet with the new quarter row ---
# A new row for 2022-Q4 is inserted at the top of the data (row 2); the
# previously-existing 2022-Q2 / 2021-Q4 rows each shift down by one row.
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

$totalSheet.Range("B4").Value = "2021-Q4"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0

$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.04

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.01
